$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new top data row for 2022-Q4, pushing
#    every existing quarter's row down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Copy the style of the last existing "A" (index) cell onto the new A7 cell
# that will be created by the row shift below.
$summary.Cells(6, 1).Copy()
$summary.Cells(7, 1).PasteSpecial(-4122)

# Shift B:D (date label / count / market value) down by one row, starting
# from the bottom so we don't clobber data we still need to read.
for ($r = 6; $r -ge 2; $r--) {
    $target = $r + 1
    $summary.Cells($target, 2).Value = $summary.Cells($r, 2).Value2
    $summary.Cells($target, 3).Value = $summary.Cells($r, 3).Value2
    $summary.Cells($target, 4).Value = $summary.Cells($r, 4).Value2
}

# Column A is a plain 0-based row index; refresh it for every data row.
for ($r = 2; $r -le 7; $r++) {
    $summary.Cells($r, 1).Value = $r - 2
}

# New row 2: the 2022-Q4 summary figures.
$summary.Cells(2, 2).Value = "2022-Q4"
$summary.Cells(2, 3).Value = 7
$summary.Cells(2, 4).Value = 0.88

# ---------------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" sheet right after "总计", holding the
#    detailed per-fund holdings for that quarter.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q4"

# Reuse formatting (header style + index-column style) from the neighbouring
# "2022-Q3" sheet so fonts/borders match the rest of the workbook.
$template = $wb.Worksheets.Item("2022-Q3")
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells(1, $i + 2).Value = $headers[$i]
}

# index, fund code, fund name, size, total position, position ratio, held value, rank
$data = @(
    @(0, "007592", "华夏价值精选混合", "4.06", "93.21", "7.49", "0.3041", 2),
    @(1, "005583", "易方达港股通红利灵活配置混合", "7.47", "88.74", "3.91", "0.2921", 8),
    @(2, "012846", "恒越蓝筹精选混合", "5.90", "86.65", "3.14", "0.1853", 5),
    @(3, "014922", "华夏ESG可持续投资一年持有混合A", "1.76", "90.74", "2.84", "0.0500", 10),
    @(4, "012993", "汇添富品牌力一年持有混合A", "1.67", "69.72", "2.27", "0.0379", 9),
    @(5, "012994", "汇添富品牌力一年持有混合C", "0.24", "69.72", "2.27", "0.0054", 9),
    @(6, "014923", "华夏ESG可持续投资一年持有混合C", "0.12", "90.74", "2.84", "0.0034", 10)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells($r, 1).Value = $row[0]
    $newSheet.Cells($r, 2).Value = "'" + $row[1]
    $newSheet.Cells($r, 3).Value = $row[2]
    $newSheet.Cells($r, 4).Value = "'" + $row[3]
    $newSheet.Cells($r, 5).Value = "'" + $row[4]
    $newSheet.Cells($r, 6).Value = "'" + $row[5]
    $newSheet.Cells($r, 7).Value = "'" + $row[6]
    $newSheet.Cells($r, 8).Value = $row[7]
    $r = $r + 1
}

Write-Output "2022-Q4 sheet inserted and 总计 sheet updated"
